$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 34021.40073141063
$ws.Range("B3").Value = 173216.2381932754
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 146842.4863074072
$ws.Range("B8").Value = 4959.039313175755
$ws.Range("B15").Value = 2
